$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.55
$ws.Range("I3").Value = 2.9
$ws.Range("L3").Value = 1.53
$ws.Range("M3").Value = 2.5
$ws.Range("N3").Value = 2.75
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 1.62
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 2.2
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 6
$ws.Range("V3").Value = 11
$ws.Range("W3").Value = 26
$ws.Range("X3").Value = 26
$ws.Range("AC3").Value = 81

# Row 4
$ws.Range("G4").Value = 1.7
$ws.Range("H4").Value = 3.8
$ws.Range("I4").Value = 4.75
$ws.Range("U4").Value = 7.5
$ws.Range("W4").Value = 13
$ws.Range("AA4").Value = 7
$ws.Range("AB4").Value = 19
$ws.Range("AD4").Value = 451
$ws.Range("AE4").Value = 11
$ws.Range("AF4").Value = 23
$ws.Range("AG4").Value = 15
$ws.Range("AJ4").Value = 41

# Row 5
$ws.Range("G5").Value = 1.6
$ws.Range("I5").Value = 5
$ws.Range("X5").Value = 12
$ws.Range("AF5").Value = 29

# Row 8
$ws.Range("N8").Value = 1.65
$ws.Range("O8").Value = 2.2
$ws.Range("R8").Value = 1.5
$ws.Range("S8").Value = 2.37
$ws.Range("AB8").Value = 12

# Row 9
$ws.Range("R9").Value = 1.58

# Row 10
$ws.Range("R10").Value = 1.47

# Row 11
$ws.Range("G11").Value = 1.7
$ws.Range("H11").Value = 3.6
$ws.Range("I11").Value = 4.2
$ws.Range("L11").Value = 1.22
$ws.Range("M11").Value = 4
$ws.Range("R11").Value = 1.63
$ws.Range("U11").Value = 9
$ws.Range("AA11").Value = 7.5
$ws.Range("AH11").Value = 51

# Row 12
$ws.Range("R12").Value = 1.63

# Row 17
$ws.Range("G17").Value = 2.35
$ws.Range("I17").Value = 3.1
$ws.Range("J17").Value = 1.06
$ws.Range("K17").Value = 10
$ws.Range("L17").Value = 1.3
$ws.Range("M17").Value = 3.4
$ws.Range("N17").Value = 2.05
$ws.Range("O17").Value = 1.75
$ws.Range("R17").Value = 1.77
$ws.Range("S17").Value = 1.92
$ws.Range("X17").Value = 21

# Row 19
$ws.Range("G19").Value = 1.7
$ws.Range("H19").Value = 3.6
$ws.Range("I19").Value = 4.2
$ws.Range("K19").Value = 7.9
$ws.Range("M19").Value = 3.55
$ws.Range("P19").Value = 1.37
$ws.Range("Q19").Value = 2.85
$ws.Range("S19").Value = 2.02
$ws.Range("AA19").Value = 7
$ws.Range("AC19").Value = 41

# Row 20
$ws.Range("G20").Value = 3.4
$ws.Range("I20").Value = 2.05
$ws.Range("N20").Value = 2
$ws.Range("O20").Value = 1.8
$ws.Range("P20").Value = 1.4
$ws.Range("Q20").Value = 2.75
$ws.Range("Y20").Value = 34
$ws.Range("Z20").Value = 9.5
$ws.Range("AF20").Value = 10
$ws.Range("AH20").Value = 19
